$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise auto-parse as a number in Excel
# (single "." => numeric). These are forced to Text before assignment,
# then restored to the default "Normal" style so no stray number format
# lingers on the cell (matches original un-styled inlineStr cells).
$numericLooking = @(
    "D5", "D7", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D17", "D19", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51"
)

foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin prices / hourly-volume percentages scraped this run,
# plus the Monero / EthereumClassic row swap (rows 27-28).
$updates = @{
    'D2' = '30.078.49'
    'E2' = '  -0.12%  '
    'D3' = '1.876.14'
    'E3' = '  -2.27%  '
    'E4' = '  +0.27%  '
    'D5' = '320.09'
    'E5' = '  -3.09%  '
    'E6' = '  +0.21%  '
    'D7' = '0.5048'
    'E7' = '  -3.19%  '
    'D8' = '0.3960'
    'E8' = '  -3.31%  '
    'D9' = '0.08215'
    'E9' = '  -3.57%  '
    'D10' = '42.13'
    'E10' = '  -1.65%  '
    'E11' = '  -3.05%  '
    'D12' = '23.53'
    'E12' = '  +4.95%  '
    'D13' = '1.867.91'
    'E13' = '  -2.16%  '
    'D14' = '6.297'
    'E14' = '  -2.14%  '
    'D15' = '7.202'
    'E15' = '  -2.84%  '
    'D16' = '1.001'
    'E16' = '  +0.09%  '
    'D17' = '91.90'
    'E17' = '  -3.89%  '
    'E18' = '  -2.42%  '
    'D19' = '0.06489'
    'E19' = '  -2.95%  '
    'E20' = '  -1.36%  '
    'D22' = '30.072.93'
    'E22' = '  -0.17%  '
    'D23' = '5.836'
    'E23' = '  -2.94%  '
    'D24' = '11.14'
    'E24' = '  -1.66%  '
    'D25' = '2.170'
    'E25' = '  -1.95%  '
    'D26' = '2.085.52'
    'E26' = '  -2.39%  '
    'B27' = 'Monero'
    'C27' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D27' = '160.70'
    'E27' = '  +0.49%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '21.15'
    'E28' = '  +0.20%  '
    'D29' = '2.241'
    'E29' = '  -8.45%  '
    'D30' = '127.56'
    'E30' = '  -1.42%  '
    'D31' = '1.081'
    'E31' = '  -0.20%  '
    'E32' = '  -1.91%  '
    'D33' = '5.952'
    'E33' = '  -1.91%  '
    'D34' = '3.698'
    'E34' = '  +1.76%  '
    'D35' = '0.02428'
    'E35' = '  -2.43%  '
    'D36' = '5.281'
    'E36' = '  +1.57%  '
    'D37' = '0.06374'
    'E37' = '  -3.86%  '
    'E38' = '  -3.73%  '
    'D39' = '1.170'
    'E39' = '  -5.37%  '
    'D40' = '8.499'
    'E40' = '  -4.54%  '
    'E41' = '  -4.01%  '
    'D42' = '1.214'
    'E42' = '  -2.81%  '
    'D43' = '11.28'
    'E43' = '  -3.11%  '
    'D44' = '13.23'
    'E44' = '  -0.52%  '
    'E45' = '  -4.17%  '
    'D46' = '2.095'
    'E46' = '  +0.41%  '
    'D47' = '3.630'
    'E47' = '  -3.78%  '
    'D48' = '122.22'
    'E48' = '  -1.98%  '
    'D49' = '1.208'
    'E49' = '  -3.42%  '
    'D50' = '77.48'
    'E50' = '  -2.83%  '
    'D51' = '1.113'
    'E51' = '  -5.45%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
